$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.948739333333333
$ws.Range("H2").Value = 14.846218
$ws.Range("I2").Value = 0.4917593264632457
$ws.Range("J2").Value = 0.4917593264632457
$ws.Range("M2").Value = 119.0164006666667
$ws.Range("N2").Value = 357.049202
$ws.Range("O2").Value = 0.9176278005170622
$ws.Range("P2").Value = 0.9176278005170622
$ws.Range("Q2").Value = 588.9811432908928
$ws.Range("R2").Value = 5300.830289618036
$ws.Range("S2").Value = 0.45125202912622
$ws.Range("T2").Value = 0.4512520291262201
$ws.Range("G3").Value = 4.948739333333333
$ws.Range("H3").Value = 14.846218
$ws.Range("I3").Value = 0.4917593264632457
$ws.Range("J3").Value = 0.4917593264632457
$ws.Range("O3").Value = 0.001755838010330732
$ws.Range("P3").Value = 0.001755838010330731
$ws.Range("Q3").Value = 1.126987955438445
$ws.Range("R3").Value = 10.142891598946
$ws.Range("S3").Value = 0.0008634497173388059
$ws.Range("T3").Value = 0.0008634497173388059
$ws.Range("G4").Value = 4.948739333333333
$ws.Range("H4").Value = 14.846218
$ws.Range("I4").Value = 0.4917593264632457
$ws.Range("J4").Value = 0.4917593264632457
$ws.Range("M4").Value = 7.816301333333333
$ws.Range("N4").Value = 23.448904
$ws.Range("O4").Value = 0.06026442877207647
$ws.Range("P4").Value = 0.06026442877207646
$ws.Range("Q4").Value = 38.68083784945244
$ws.Range("R4").Value = 348.127540645072
$ws.Range("S4").Value = 0.02963559490264856
$ws.Range("T4").Value = 0.02963559490264856
$ws.Range("G5").Value = 4.948739333333333
$ws.Range("H5").Value = 14.846218
$ws.Range("I5").Value = 0.4917593264632457
$ws.Range("J5").Value = 0.4917593264632457
$ws.Range("M5").Value = 0.105045
$ws.Range("N5").Value = 0.315135
$ws.Range("O5").Value = 0.000809906968832672
$ws.Range("P5").Value = 0.000809906968832672
$ws.Range("Q5").Value = 0.51984032327
$ws.Range("R5").Value = 4.67856290943
$ws.Range("S5").Value = 0.0003982793054910437
$ws.Range("T5").Value = 0.0003982793054910437
$ws.Range("G6").Value = 4.948739333333333
$ws.Range("H6").Value = 14.846218
$ws.Range("I6").Value = 0.4917593264632457
$ws.Range("J6").Value = 0.4917593264632457
$ws.Range("M6").Value = 2.534602333333333
$ws.Range("N6").Value = 7.603807
$ws.Range("O6").Value = 0.01954202573169801
$ws.Range("P6").Value = 0.01954202573169801
$ws.Range("Q6").Value = 12.54308626132511
$ws.Range("R6").Value = 112.887776351926
$ws.Range("S6").Value = 0.009609973411547229
$ws.Range("T6").Value = 0.009609973411547231
$ws.Range("G7").Value = 3.979395333333333
$ws.Range("I7").Value = 0.395435006178203
$ws.Range("J7").Value = 0.395435006178203
$ws.Range("M7").Value = 119.0164006666667
$ws.Range("N7").Value = 357.049202
$ws.Range("O7").Value = 0.9176278005170622
$ws.Range("P7").Value = 0.9176278005170622
$ws.Range("Q7").Value = 473.6133094030635
$ws.Range("R7").Value = 4262.519784627571
$ws.Range("S7").Value = 0.3628621549667553
$ws.Range("T7").Value = 0.3628621549667553
$ws.Range("G8").Value = 3.979395333333333
$ws.Range("I8").Value = 0.395435006178203
$ws.Range("J8").Value = 0.395435006178203
$ws.Range("O8").Value = 0.001755838010330732
$ws.Range("P8").Value = 0.001755838010330731
$ws.Range("Q8").Value = 0.9062369845157777
$ws.Range("R8").Value = 8.156132860642
$ws.Range("S8").Value = 0.0006943198144630565
$ws.Range("T8").Value = 0.0006943198144630565
$ws.Range("G9").Value = 3.979395333333333
$ws.Range("I9").Value = 0.395435006178203
$ws.Range("J9").Value = 0.395435006178203
$ws.Range("M9").Value = 7.816301333333333
$ws.Range("N9").Value = 23.448904
$ws.Range("O9").Value = 0.06026442877207647
$ws.Range("P9").Value = 0.06026442877207646
$ws.Range("Q9").Value = 31.10415304979377
$ws.Range("R9").Value = 279.937377448144
$ws.Range("S9").Value = 0.02383066476381193
$ws.Range("T9").Value = 0.02383066476381193
$ws.Range("G10").Value = 3.979395333333333
$ws.Range("I10").Value = 0.395435006178203
$ws.Range("J10").Value = 0.395435006178203
$ws.Range("M10").Value = 0.105045
$ws.Range("N10").Value = 0.315135
$ws.Range("O10").Value = 0.000809906968832672
$ws.Range("P10").Value = 0.000809906968832672
$ws.Range("Q10").Value = 0.41801558279
$ws.Range("R10").Value = 3.76214024511
$ws.Range("S10").Value = 0.0003202655672241173
$ws.Range("T10").Value = 0.0003202655672241173
$ws.Range("G11").Value = 3.979395333333333
$ws.Range("I11").Value = 0.395435006178203
$ws.Range("J11").Value = 0.395435006178203
$ws.Range("M11").Value = 2.534602333333333
$ws.Range("N11").Value = 7.603807
$ws.Range("O11").Value = 0.01954202573169801
$ws.Range("P11").Value = 0.01954202573169801
$ws.Range("Q11").Value = 10.08618469712244
$ws.Range("R11").Value = 90.775662274102
$ws.Range("S11").Value = 0.007727601065948605
$ws.Range("T11").Value = 0.007727601065948607
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.113241666666667
$ws.Range("H12").Value = 3.339725
$ws.Range("I12").Value = 0.1106235215306998
$ws.Range("J12").Value = 0.1106235215306998
$ws.Range("M12").Value = 119.0164006666667
$ws.Range("N12").Value = 357.049202
$ws.Range("O12").Value = 0.9176278005170622
$ws.Range("P12").Value = 0.9176278005170622
$ws.Range("Q12").Value = 132.4940162388278
$ws.Range("R12").Value = 1192.44614614945
$ws.Range("S12").Value = 0.101511218747668
$ws.Range("T12").Value = 0.101511218747668
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.113241666666667
$ws.Range("H13").Value = 3.339725
$ws.Range("I13").Value = 0.1106235215306998
$ws.Range("J13").Value = 0.1106235215306998
$ws.Range("O13").Value = 0.001755838010330732
$ws.Range("P13").Value = 0.001755838010330731
$ws.Range("Q13").Value = 0.253521122313889
$ws.Range("R13").Value = 2.281690100825001
$ws.Range("S13").Value = 0.0001942369839402428
$ws.Range("T13").Value = 0.0001942369839402428
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.113241666666667
$ws.Range("H14").Value = 3.339725
$ws.Range("I14").Value = 0.1106235215306998
$ws.Range("J14").Value = 0.1106235215306998
$ws.Range("M14").Value = 7.816301333333333
$ws.Range("N14").Value = 23.448904
$ws.Range("O14").Value = 0.06026442877207647
$ws.Range("P14").Value = 0.06026442877207646
$ws.Range("Q14").Value = 8.70143232348889
$ws.Range("R14").Value = 78.3128909114
$ws.Range("S14").Value = 0.006666663333803127
$ws.Range("T14").Value = 0.006666663333803127
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.113241666666667
$ws.Range("H15").Value = 3.339725
$ws.Range("I15").Value = 0.1106235215306998
$ws.Range("J15").Value = 0.1106235215306998
$ws.Range("M15").Value = 0.105045
$ws.Range("N15").Value = 0.315135
$ws.Range("O15").Value = 0.000809906968832672
$ws.Range("P15").Value = 0.000809906968832672
$ws.Range("Q15").Value = 0.116940470875
$ws.Range("R15").Value = 1.052464237875
$ws.Range("S15").Value = 0.00008959476100452492
$ws.Range("T15").Value = 0.00008959476100452492
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.113241666666667
$ws.Range("H16").Value = 3.339725
$ws.Range("I16").Value = 0.1106235215306998
$ws.Range("J16").Value = 0.1106235215306998
$ws.Range("M16").Value = 2.534602333333333
$ws.Range("N16").Value = 7.603807
$ws.Range("O16").Value = 0.01954202573169801
$ws.Range("P16").Value = 0.01954202573169801
$ws.Range("Q16").Value = 2.821624925897223
$ws.Range("R16").Value = 25.394624333075
$ws.Range("S16").Value = 0.002161807704283985
$ws.Range("T16").Value = 0.002161807704283985
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.02195966666666667
$ws.Range("H17").Value = 0.06587899999999999
$ws.Range("I17").Value = 0.002182145827851387
$ws.Range("J17").Value = 0.002182145827851387
$ws.Range("M17").Value = 119.0164006666667
$ws.Range("N17").Value = 357.049202
$ws.Range("O17").Value = 0.9176278005170622
$ws.Range("P17").Value = 0.9176278005170622
$ws.Range("Q17").Value = 2.613560486506444
$ws.Range("R17").Value = 23.522044378558
$ws.Range("S17").Value = 0.002002397676418752
$ws.Range("T17").Value = 0.002002397676418752
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.02195966666666667
$ws.Range("H18").Value = 0.06587899999999999
$ws.Range("I18").Value = 0.002182145827851387
$ws.Range("J18").Value = 0.002182145827851387
$ws.Range("O18").Value = 0.001755838010330732
$ws.Range("P18").Value = 0.001755838010330731
$ws.Range("Q18").Value = 0.005000926129222222
$ws.Range("R18").Value = 0.045008335163
$ws.Range("S18").Value = 0.000003831494588626086
$ws.Range("T18").Value = 0.000003831494588626085
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.02195966666666667
$ws.Range("H19").Value = 0.06587899999999999
$ws.Range("I19").Value = 0.002182145827851387
$ws.Range("J19").Value = 0.002182145827851387
$ws.Range("M19").Value = 7.816301333333333
$ws.Range("N19").Value = 23.448904
$ws.Range("O19").Value = 0.06026442877207647
$ws.Range("P19").Value = 0.06026442877207646
$ws.Range("Q19").Value = 0.1716433718462222
$ws.Range("R19").Value = 1.544790346616
$ws.Range("S19").Value = 0.0001315057718128337
$ws.Range("T19").Value = 0.0001315057718128337
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0.3333333333333333
$ws.Range("G20").Value = 0.02195966666666667
$ws.Range("H20").Value = 0.06587899999999999
$ws.Range("I20").Value = 0.002182145827851387
$ws.Range("J20").Value = 0.002182145827851387
$ws.Range("M20").Value = 0.105045
$ws.Range("N20").Value = 0.315135
$ws.Range("O20").Value = 0.000809906968832672
$ws.Range("P20").Value = 0.000809906968832672
$ws.Range("Q20").Value = 0.002306753185
$ws.Range("R20").Value = 0.020760778665
$ws.Range("S20").Value = 0.000001767335112985978
$ws.Range("T20").Value = 0.000001767335112985978
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 0.3333333333333333
$ws.Range("G21").Value = 0.02195966666666667
$ws.Range("H21").Value = 0.06587899999999999
$ws.Range("I21").Value = 0.002182145827851387
$ws.Range("J21").Value = 0.002182145827851387
$ws.Range("M21").Value = 2.534602333333333
$ws.Range("N21").Value = 7.603807
$ws.Range("O21").Value = 0.01954202573169801
$ws.Range("P21").Value = 0.01954202573169801
$ws.Range("Q21").Value = 0.05565902237255555
$ws.Range("R21").Value = 0.500931201353
$ws.Range("S21").Value = 0.00004264354991818926
$ws.Range("T21").Value = 0.00004264354991818926
